$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Remove the duplicate empty spacer paragraph near the end of the document.
#    Before: ... "Zijn de classes..." / ind=360 / ind=720 / (none) / ind=360 / ind=360
#    After  (content-wise, before the new paragraphs are added): the duplicate
#    "ind=360" paragraph right before the very last paragraph is removed.
#    At this point (before any insertions) this is paragraph index 20.
# ---------------------------------------------------------------------------
$dupPara = $d.Paragraphs.Item(20)
$dupPara.Range.Delete()

# ---------------------------------------------------------------------------
# 2. Locate the "Zijn de classes nu te overlappend?" paragraph (still #16,
#    unaffected by the deletion above) and drop the _GoBack bookmark that
#    currently sits on it -- it will be re-created on the new last question.
# ---------------------------------------------------------------------------
$lastQuestion = $d.Paragraphs.Item(16)

if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# 3. Insert the four new paragraphs as a single chain right after
#    "Zijn de classes...". Inserting them one after another lets each new
#    paragraph naturally inherit the "NoSpacing" style and numId=2 bullet
#    list from its predecessor, so the final paragraph keeps reusing the
#    existing list (numId 2) instead of Word fabricating a brand new list.
# ---------------------------------------------------------------------------
$lastQuestion.Range.InsertParagraphAfter()
$spacer1 = $d.Paragraphs.Item(17)

$spacer1.Range.InsertParagraphAfter()
$spacer2 = $d.Paragraphs.Item(18)

$spacer2.Range.InsertParagraphAfter()
$heading3 = $d.Paragraphs.Item(19)

$heading3.Range.InsertParagraphAfter()
$newQuestion = $d.Paragraphs.Item(20)

# Spacer paragraphs: plain "No Spacing" paragraphs indented like the other
# spacer paragraphs already present in the document, and not part of a list.
$spacer1.Range.ListFormat.RemoveNumbers()
$spacer1.Range.ParagraphFormat.LeftIndent = 18   # 360 twips

$spacer2.Range.ListFormat.RemoveNumbers()
$spacer2.Range.ParagraphFormat.LeftIndent = 36   # 720 twips

# "Vragen 3e keer:" heading paragraph, not part of the bullet list.
$heading3.Range.ListFormat.RemoveNumbers()
$heading3.Range.Text = "Vragen 3e keer:"
$eStart = $heading3.Range.Start + 8
$eRange = $d.Range($eStart, $eStart + 1)
$eRange.Font.Superscript = $true

# New bullet question, keeps numId=2 (inherited through the insert chain).
$newQuestion.Range.Text = "Github: Moet er in iedere map een nieuwe READme met uitleg, net als in de voorbeeldrepo?"

$bmStart = $newQuestion.Range.Start
$bmEnd = $newQuestion.Range.End - 1
$bmRange = $d.Range($bmStart, $bmEnd)
$d.Bookmarks.Add("_GoBack", $bmRange)
